$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 1 (headers) - values stay the same, but a new blank trailing
# cell (S1) is added which extends the used range to column S.
# ---------------------------------------------------------------
$ws.Range("S1").Font.Name = "Arial"
$ws.Range("S1").Font.Size = 10

# ---------------------------------------------------------------
# Row 2 - update existing values to the new data set.
# ---------------------------------------------------------------
$ws.Range("A2").Value = 123456
$ws.Range("B2").Value = "Sada123"
$ws.Range("D2").Value = 218884411
$ws.Range("E2").Value = "sada@gmail.com"
$ws.Range("F2").Value = "WITA"
$ws.Range("G2").Value = 12141
$ws.Range("H2").Value = "BCA"
$ws.Range("J2").Value = "baru1"
$ws.Range("K2").Value = "Laki-laki"
$ws.Range("N2").Value = "Baru1, baru2"
$ws.Range("O2").Value = "SD"
$ws.Range("P2").Value = "subarea456"
$ws.Range("Q2").Value = "area124"
$ws.Range("R2").Value = "region123"
$ws.Range("S2").Font.Name = "Arial"
$ws.Range("S2").Font.Size = 10

# ---------------------------------------------------------------
# Row 3 - brand new row of data.
# ---------------------------------------------------------------
$ws.Range("A3").Value = 12345
$ws.Range("B3").Value = "Sada12345"
$ws.Range("D3").Value = 8889996
$ws.Range("E3").Value = "sada123@gmail.com"
$ws.Range("F3").Value = "WIB"
$ws.Range("J3").Value = "baru2"
$ws.Range("K3").Value = "Perempuan"
$ws.Range("N3").Value = "Baru4, baru5"
$ws.Range("O3").Value = "SD"
$ws.Range("P3").Value = "subarea456"
$ws.Range("Q3").Value = "area1245"
$ws.Range("R3").Value = "region123"

# Row 3 uses a 15.75pt custom height, same as the other rows.
$ws.Rows.Item(3).RowHeight = 15.75

# ---------------------------------------------------------------
# Hyperlinks - rebuild the mailto: hyperlinks on E2 and E3 so the
# relationship parts point at the right addresses.
# ---------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:sada@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:sada123@gmail.com")

# ---------------------------------------------------------------
# Column widths for the newly used E/F columns.
# ---------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 18.5
$ws.Columns.Item(6).ColumnWidth = 24.3

# ---------------------------------------------------------------
# View state - move the selection to R3 (matches the saved
# selection in the authored workbook) and scroll so column G is
# the first visible column.
# ---------------------------------------------------------------
$ws.Range("R3").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
